$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay text (matches source formatting)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "64.211.03"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").Value = "2.743.15"
$ws.Range("E3").Value = "  -0.52%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "573.30"
$ws.Range("E5").Value = "  -1.15%  "

# Row 6
$ws.Range("D6").Value = "158.66"
$ws.Range("E6").Value = "  -1.15%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("D8").Value = "0.600"
$ws.Range("E8").Value = "  -2.13%  "

# Row 9
$ws.Range("E9").Value = "  -2.61%  "

# Row 10
$ws.Range("D10").Value = "0.166"
$ws.Range("E10").Value = "  +4.79%  "

# Row 11
$ws.Range("D11").Value = "5.77"
$ws.Range("E11").Value = "  -1.62%  "

# Row 12
$ws.Range("D12").Value = "0.384"
$ws.Range("E12").Value = "  -2.45%  "

# Row 13
$ws.Range("D13").Value = "3.230.81"
$ws.Range("E13").Value = "  -0.31%  "

# Row 14
$ws.Range("D14").Value = "26.97"
$ws.Range("E14").Value = "  -0.57%  "

# Row 15
$ws.Range("D15").Value = "63.903.41"
$ws.Range("E15").Value = "  -0.03%  "

# Row 16
$ws.Range("D16").Value = "0.0000151"
$ws.Range("E16").Value = "  -3.12%  "

# Row 17
$ws.Range("D17").Value = "2.750.28"
$ws.Range("E17").Value = "  -0.23%  "

# Row 18
$ws.Range("D18").Value = "12.10"
$ws.Range("E18").Value = "  -1.95%  "

# Row 19
$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  -3.19%  "

# Row 20
$ws.Range("D20").Value = "355.84"
$ws.Range("E20").Value = "  -2.11%  "

# Row 21
$ws.Range("D21").Value = "6.62"
$ws.Range("E21").Value = "  -4.19%  "

# Row 22
$ws.Range("E22").Value = "  -0.09%  "

# Row 23
$ws.Range("D23").Value = "0.523"
$ws.Range("E23").Value = "  -8.42%  "

# Row 24
$ws.Range("D24").Value = "64.79"
$ws.Range("E24").Value = "  -2.53%  "

# Row 25
$ws.Range("D25").Value = "0.170"
$ws.Range("E25").Value = "  -1.56%  "

# Row 26
$ws.Range("D26").Value = "8.51"
$ws.Range("E26").Value = "  -2.10%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0912"
$ws.Range("E28").Value = "  -3.46%  "

# Row 29
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  +2.46%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.96"
$ws.Range("E30").Value = "  -2.06%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").Value = "1.36"
$ws.Range("E31").Value = "  +7.67%  "

# Row 32
$ws.Range("D32").Value = "166.46"
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "1.50"
$ws.Range("E33").Value = "  +0.95%  "

# Row 34
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "20.10"
$ws.Range("E34").Value = "  -2.55%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "4.91"
$ws.Range("E35").Value = "  -2.57%  "

# Row 36
$ws.Range("E36").Value = "  +0.00%  "

# Row 37
$ws.Range("E37").Value = "  +0.43%  "

# Row 38
$ws.Range("D38").Value = "0.995"
$ws.Range("E38").Value = "  -2.42%  "

# Row 39
$ws.Range("D39").Value = "348.33"
$ws.Range("E39").Value = "  +2.97%  "

# Row 40
$ws.Range("D40").Value = "6.34"
$ws.Range("E40").Value = "  +2.42%  "

# Row 41
$ws.Range("D41").Value = "4.16"
$ws.Range("E41").Value = "  -1.87%  "

# Row 42
$ws.Range("D42").Value = "39.04"
$ws.Range("E42").Value = "  -1.28%  "

# Row 43
$ws.Range("D43").Value = "22.35"
$ws.Range("E43").Value = "  +1.05%  "

# Row 44
$ws.Range("D44").Value = "21.41"
$ws.Range("E44").Value = "  -3.72%  "

# Row 45
$ws.Range("D45").Value = "0.0589"
$ws.Range("E45").Value = "  -2.25%  "

# Row 46
$ws.Range("D46").Value = "135.44"
$ws.Range("E46").Value = "  -1.03%  "

# Row 47
$ws.Range("D47").Value = "0.626"
$ws.Range("E47").Value = "  -2.35%  "

# Row 48
$ws.Range("E48").Value = "  -1.62%  "

# Row 49
$ws.Range("D49").Value = "0.0251"
$ws.Range("E49").Value = "  -3.31%  "

# Row 50
$ws.Range("D50").Value = "0.998"
$ws.Range("E50").Value = "  -0.03%  "

# Row 51
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "11.02"
$ws.Range("E51").Value = "  -0.22%  "
